# Auto-update draw results: append the 2025-10-07 Pick 4 draw as row 21.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A21:E21")

# The sheet stores every value (including number-looking ones like the
# phase code "251007") as plain text, so force Text format before writing
# the values to stop Excel from auto-coercing them into dates/numbers.
$newRow.NumberFormat = "@"

$ws.Range("A21").Value = "2025-10-07"
$ws.Range("B21").Value = "Pick 4"
$ws.Range("C21").Value = "251007"
$ws.Range("D21").Value = "2-4-2-6"
$ws.Range("E21").Value = "2025-10-07T21:37:50.240+04:00"

# Restore the default cell style so the new row matches the formatting
# (no explicit style) used by the rest of the table.
$newRow.Style = "Normal"
